# Edit script: apply 2025-09 Bundesliga matchday-3 fixtures insertion and
# re-align subsequent match rows (data pipeline shift), extending the table
# from A1:L696 to A1:L705.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(540425, '2025-09-12T18:30:00Z', 3, 3, 'Bayer 04 Leverkusen', 19, 'Eintracht Frankfurt', 3, 1, 'HomeWin', 4, 6),
    @(540426, '2025-09-13T13:30:00Z', 3, 17, 'SC Freiburg', 10, 'VfB Stuttgart', 3, 1, 'HomeWin', 6, 3),
    @(540427, '2025-09-13T13:30:00Z', 3, 15, '1. FSV Mainz 05', 721, 'RB Leipzig', 0, 1, 'AwayWin', 3, 4),
    @(540429, '2025-09-13T13:30:00Z', 3, 11, 'VfL Wolfsburg', 1, '1. FC Köln', 3, 3, 'Draw', 3, 4),
    @(540430, '2025-09-13T13:30:00Z', 3, 28, '1. FC Union Berlin', 2, 'TSG 1899 Hoffenheim', 2, 4, 'AwayWin', 0, 7),
    @(540432, '2025-09-13T13:30:00Z', 3, 44, '1. FC Heidenheim 1846', 4, 'Borussia Dortmund', 0, 2, 'AwayWin', 0, 3),
    @(540424, '2025-09-13T16:30:00Z', 3, 5, 'FC Bayern München', 7, 'Hamburger SV', 5, 0, 'HomeWin', 0, 1),
    @(540431, '2025-09-14T13:30:00Z', 3, 20, 'FC St. Pauli', 16, 'FC Augsburg', 2, 1, 'HomeWin', 4, 4),
    @(540428, '2025-09-14T15:30:00Z', 3, 18, 'Borussia Mönchengladbach', 12, 'SV Werder Bremen', 0, 4, 'AwayWin', 1, 0),
    @(540433, '2024-08-02T00:00:00Z', 0, 722, '1. FC Köln', 1722, 'Hamburger SV', 1, 2, 'A', 0, 0),
    @(540434, '2024-08-10T00:00:00Z', 0, 723, 'Elversberg', 1723, '1. FC Köln', 2, 2, 'D', 4, 0),
    @(540435, '2024-08-10T00:00:00Z', 0, 724, 'Hamburger SV', 1724, 'Hertha', 1, 1, 'D', 0, 0),
    @(540436, '2024-08-23T00:00:00Z', 0, 725, 'Hannover', 1725, 'Hamburger SV', 1, 0, 'H', 0, 0),
    @(540437, '2024-08-24T00:00:00Z', 0, 726, '1. FC Köln', 1726, 'Braunschweig', 5, 0, 'H', 3, 0),
    @(540438, '2024-08-31T00:00:00Z', 0, 727, 'Hamburger SV', 1727, 'Preußen Münster', 4, 1, 'H', 0, 0),
    @(540439, '2024-09-01T00:00:00Z', 0, 728, 'Schalke 04', 1728, '1. FC Köln', 1, 3, 'A', 0, 0),
    @(540440, '2024-09-14T00:00:00Z', 0, 729, '1. FC Köln', 1729, 'Magdeburg', 1, 2, 'A', 0, 0),
    @(540441, '2024-09-15T00:00:00Z', 0, 730, 'Hamburger SV', 1730, 'Regensburg', 5, 0, 'H', 0, 0),
    @(540442, '2024-09-21T00:00:00Z', 0, 731, 'Fortuna Dusseldorf', 1731, '1. FC Köln', 2, 2, 'D', 0, 0),
    @(540443, '2024-09-21T00:00:00Z', 0, 732, 'Kaiserslautern', 1732, 'Hamburger SV', 2, 2, 'D', 0, 0),
    @(540444, '2024-09-28T00:00:00Z', 0, 733, 'Hamburger SV', 1733, 'Paderborn', 2, 2, 'D', 0, 0),
    @(540445, '2024-09-29T00:00:00Z', 0, 734, '1. FC Köln', 1734, 'Karlsruhe', 4, 4, 'D', 0, 0),
    @(540446, '2024-10-05T00:00:00Z', 0, 735, '1. FC Köln', 1735, 'Ulm', 2, 0, 'H', 0, 0),
    @(540447, '2024-10-06T00:00:00Z', 0, 736, 'Fortuna Dusseldorf', 1736, 'Hamburger SV', 0, 3, 'A', 0, 0),
    @(540448, '2024-10-18T00:00:00Z', 0, 737, 'Darmstadt', 1737, '1. FC Köln', 5, 1, 'H', 0, 0),
    @(540449, '2024-10-20T00:00:00Z', 0, 738, 'Hamburger SV', 1738, 'Magdeburg', 3, 1, 'H', 0, 0),
    @(540450, '2024-10-25T00:00:00Z', 0, 739, '1. FC Köln', 1739, 'Paderborn', 1, 2, 'A', 0, 0),
    @(540451, '2024-10-26T00:00:00Z', 0, 740, 'Elversberg', 1740, 'Hamburger SV', 4, 2, 'H', 0, 0),
    @(540452, '2024-11-02T00:00:00Z', 0, 741, 'Hertha', 1741, '1. FC Köln', 0, 1, 'A', 0, 0),
    @(540453, '2024-11-03T00:00:00Z', 0, 742, 'Hamburger SV', 1742, 'Nurnberg', 1, 1, 'D', 0, 0),
    @(540454, '2024-11-08T00:00:00Z', 0, 743, 'Braunschweig', 1743, 'Hamburger SV', 3, 1, 'H', 0, 0),
    @(540455, '2024-11-09T00:00:00Z', 0, 744, '1. FC Köln', 1744, 'Greuther Furth', 1, 0, 'H', 0, 0),
    @(540456, '2024-11-22T00:00:00Z', 0, 745, 'Preußen Münster', 1745, '1. FC Köln', 0, 1, 'A', 0, 0),
    @(540457, '2024-11-23T00:00:00Z', 0, 746, 'Hamburger SV', 1746, 'Schalke 04', 2, 2, 'D', 0, 0),
    @(540458, '2024-11-30T00:00:00Z', 0, 747, '1. FC Köln', 1747, 'Hannover', 2, 2, 'D', 0, 0),
    @(540459, '2024-12-01T00:00:00Z', 0, 748, 'Karlsruhe', 1748, 'Hamburger SV', 1, 3, 'A', 0, 0),
    @(540460, '2024-12-08T00:00:00Z', 0, 749, 'Hamburger SV', 1749, 'Darmstadt', 2, 2, 'D', 0, 0),
    @(540461, '2024-12-08T00:00:00Z', 0, 750, 'Regensburg', 1750, '1. FC Köln', 0, 1, 'A', 0, 0),
    @(540462, '2024-12-14T00:00:00Z', 0, 751, 'Ulm', 1751, 'Hamburger SV', 1, 1, 'D', 0, 0),
    @(540463, '2024-12-15T00:00:00Z', 0, 752, '1. FC Köln', 1752, 'Nurnberg', 3, 1, 'H', 0, 0),
    @(540464, '2024-12-21T00:00:00Z', 0, 753, 'Hamburger SV', 1753, 'Greuther Furth', 5, 0, 'H', 0, 0),
    @(540465, '2024-12-22T00:00:00Z', 0, 754, 'Kaiserslautern', 1754, '1. FC Köln', 0, 1, 'A', 0, 0),
    @(540466, '2025-01-18T00:00:00Z', 0, 755, 'Hamburger SV', 1755, '1. FC Köln', 1, 0, 'H', 0, 0),
    @(540467, '2025-01-25T00:00:00Z', 0, 756, '1. FC Köln', 1756, 'Elversberg', 1, 0, 'H', 0, 0),
    @(540468, '2025-01-25T00:00:00Z', 0, 757, 'Hertha', 1757, 'Hamburger SV', 2, 3, 'A', 0, 0),
    @(540469, '2025-02-01T00:00:00Z', 0, 758, 'Braunschweig', 1758, '1. FC Köln', 1, 2, 'A', 0, 0),
    @(540470, '2025-02-02T00:00:00Z', 0, 759, 'Hamburger SV', 1759, 'Hannover', 2, 2, 'D', 0, 0),
    @(540471, '2025-02-07T00:00:00Z', 0, 760, 'Preußen Münster', 1760, 'Hamburger SV', 1, 2, 'A', 0, 0),
    @(540472, '2025-02-09T00:00:00Z', 0, 761, '1. FC Köln', 1761, 'Schalke 04', 1, 0, 'H', 0, 0),
    @(540473, '2025-02-14T00:00:00Z', 0, 762, 'Magdeburg', 1762, '1. FC Köln', 3, 0, 'H', 0, 0),
    @(540474, '2025-02-16T00:00:00Z', 0, 763, 'Regensburg', 1763, 'Hamburger SV', 1, 1, 'D', 0, 0),
    @(540475, '2025-02-21T00:00:00Z', 0, 764, 'Hamburger SV', 1764, 'Kaiserslautern', 3, 0, 'H', 0, 0),
    @(540476, '2025-02-23T00:00:00Z', 0, 765, '1. FC Köln', 1765, 'Fortuna Dusseldorf', 1, 1, 'D', 0, 0),
    @(540477, '2025-03-01T00:00:00Z', 0, 766, 'Karlsruhe', 1766, '1. FC Köln', 1, 0, 'H', 0, 0),
    @(540478, '2025-03-02T00:00:00Z', 0, 767, 'Paderborn', 1767, 'Hamburger SV', 2, 0, 'H', 0, 0),
    @(540479, '2025-03-08T00:00:00Z', 0, 768, 'Ulm', 1768, '1. FC Köln', 0, 1, 'A', 0, 0),
    @(540480, '2025-03-08T00:00:00Z', 0, 769, 'Hamburger SV', 1769, 'Fortuna Dusseldorf', 4, 1, 'H', 0, 0),
    @(540481, '2025-03-14T00:00:00Z', 0, 770, 'Magdeburg', 1770, 'Hamburger SV', 0, 3, 'A', 0, 0),
    @(540482, '2025-03-15T00:00:00Z', 0, 771, '1. FC Köln', 1771, 'Darmstadt', 2, 1, 'H', 0, 0),
    @(540483, '2025-03-28T00:00:00Z', 0, 772, 'Hamburger SV', 1772, 'Elversberg', 0, 0, 'D', 0, 0),
    @(540484, '2025-03-29T00:00:00Z', 0, 773, 'Paderborn', 1773, '1. FC Köln', 1, 2, 'A', 0, 0),
    @(540485, '2025-04-05T00:00:00Z', 0, 774, 'Nurnberg', 1774, 'Hamburger SV', 0, 3, 'A', 0, 0),
    @(540486, '2025-04-05T00:00:00Z', 0, 775, '1. FC Köln', 1775, 'Hertha', 0, 1, 'A', 0, 0),
    @(540487, '2025-04-11T00:00:00Z', 0, 776, 'Greuther Furth', 1776, '1. FC Köln', 1, 1, 'D', 0, 0),
    @(540488, '2025-04-11T00:00:00Z', 0, 777, 'Hamburger SV', 1777, 'Braunschweig', 2, 4, 'A', 0, 0),
    @(540489, '2025-04-19T00:00:00Z', 0, 778, 'Schalke 04', 1778, 'Hamburger SV', 2, 2, 'D', 0, 0),
    @(540490, '2025-04-20T00:00:00Z', 0, 779, '1. FC Köln', 1779, 'Preußen Münster', 3, 1, 'H', 0, 0),
    @(540491, '2025-04-27T00:00:00Z', 0, 780, 'Hamburger SV', 1780, 'Karlsruhe', 1, 2, 'A', 0, 0),
    @(540492, '2025-04-27T00:00:00Z', 0, 781, 'Hannover', 1781, '1. FC Köln', 1, 0, 'H', 0, 0),
    @(540493, '2025-05-03T00:00:00Z', 0, 782, 'Darmstadt', 1782, 'Hamburger SV', 0, 4, 'A', 0, 0),
    @(540494, '2025-05-03T00:00:00Z', 0, 783, '1. FC Köln', 1783, 'Regensburg', 1, 1, 'D', 0, 0),
    @(540495, '2025-05-09T00:00:00Z', 0, 784, 'Nurnberg', 1784, '1. FC Köln', 1, 2, 'A', 0, 0),
    @(540496, '2025-05-10T00:00:00Z', 0, 785, 'Hamburger SV', 1785, 'Ulm', 6, 1, 'H', 0, 0),
    @(540497, '2025-05-18T00:00:00Z', 0, 786, '1. FC Köln', 1786, 'Kaiserslautern', 4, 0, 'H', 0, 0),
    @(540498, '2025-05-18T00:00:00Z', 0, 787, 'Greuther Furth', 1787, 'Hamburger SV', 3, 2, 'H', 0, 0)
)

$startRow = 631
for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $r = $startRow + $idx
    $row = $data[$idx]
    $ws.Cells.Item($r, 1).Value  = $row[0]
    $ws.Cells.Item($r, 2).Value  = $row[1]
    $ws.Cells.Item($r, 3).Value  = $row[2]
    $ws.Cells.Item($r, 4).Value  = $row[3]
    $ws.Cells.Item($r, 5).Value  = $row[4]
    $ws.Cells.Item($r, 6).Value  = $row[5]
    $ws.Cells.Item($r, 7).Value  = $row[6]
    $ws.Cells.Item($r, 8).Value  = $row[7]
    $ws.Cells.Item($r, 9).Value  = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
}

Write-Output "Rows 631-705 written."
